# Updates the "Fecha" (D), "Volumen" (J), "Precio mínimo" (K),
# "Precio máximo" (L), "Precio promedio ponderado" (M) and "Precio $/Kg" (P)
# columns for rows 2-12 with the new weekly values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (Date, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg)
$rows = @{
    2  = @(44592, 120, 12000, 13000, 12500, 962)
    3  = @(44229, 120, 44000, 45000, 44500, 3423)
    4  = @(44159, 100, 23000, 24000, 23500, 1808)
    5  = @(44379, 120, 12000, 13000, 12667, 974)
    6  = @(44469, 140, 13000, 14000, 13500, 1038)
    7  = @(44320, 160, 19000, 20000, 19500, 1500)
    8  = @(44616, 120, 19000, 20000, 19500, 1500)
    9  = @(44406, 160, 17000, 18000, 17500, 1346)
    10 = @(44580, 160, 11000, 12000, 11500, 885)
    11 = @(44389, 120, 12000, 13000, 12500, 962)
    12 = @(44397, 140, 12500, 13000, 12750, 981)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("J$r").Value = $vals[1]
    $ws.Range("K$r").Value = $vals[2]
    $ws.Range("L$r").Value = $vals[3]
    $ws.Range("M$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
}
